# Insert 2 new data rows at the top of the "Apio" data block (before the
# existing row 471), pushing all subsequent rows (old 471-568) down to
# 473-570. Then populate the two newly inserted rows with fresh data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 471; existing rows shift down by 2.
$ws.Rows("471:472").Insert()

# --- New row 471 ---
$ws.Cells.Item(471, 1).Value = 10
$ws.Cells.Item(471, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(471, 3).Value = "La Araucanía"
$ws.Cells.Item(471, 4).Value = 45209
$ws.Cells.Item(471, 5).Value = 9
$ws.Cells.Item(471, 6).Value = 100112017
$ws.Cells.Item(471, 7).Value = "Apio"
$ws.Cells.Item(471, 8).Value = "Americana (o)"
$ws.Cells.Item(471, 9).Value = "Primera"
$ws.Cells.Item(471, 10).Value = 300
$ws.Cells.Item(471, 11).Value = 8000
$ws.Cells.Item(471, 12).Value = 8000
$ws.Cells.Item(471, 13).Value = 8000
$ws.Cells.Item(471, 14).Value = "$/caja 8 unidades"
$ws.Cells.Item(471, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(471, 16).Value = 8000
$ws.Cells.Item(471, 17).Value = 1
$ws.Cells.Item(471, 18).Value = "Hortaliza"

# --- New row 472 ---
$ws.Cells.Item(472, 1).Value = 10
$ws.Cells.Item(472, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(472, 3).Value = "La Araucanía"
$ws.Cells.Item(472, 4).Value = 45209
$ws.Cells.Item(472, 5).Value = 9
$ws.Cells.Item(472, 6).Value = 100112017
$ws.Cells.Item(472, 7).Value = "Apio"
$ws.Cells.Item(472, 8).Value = "Americana (o)"
$ws.Cells.Item(472, 9).Value = "Primera"
$ws.Cells.Item(472, 10).Value = 110
$ws.Cells.Item(472, 11).Value = 7000
$ws.Cells.Item(472, 12).Value = 8000
$ws.Cells.Item(472, 13).Value = 7545
$ws.Cells.Item(472, 14).Value = "$/docena de matas"
$ws.Cells.Item(472, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(472, 16).Value = 1258
$ws.Cells.Item(472, 17).Value = 6
$ws.Cells.Item(472, 18).Value = "Hortaliza"
